# Add the new "168_antibiotic_11.19" plate-layout sheet after the last
# existing sheet (mgz1_antibiotic_11.9), populate its header + data rows,
# match the look (bold header row, column C width, 140% zoom) of the
# sibling plate-layout sheets, and leave it as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "168_antibiotic_11.19"

# --- header row ---------------------------------------------------------
$headers = @("well", "strain", "compound", "concentration_unit", "concentration", "compound_volume", "media", "media_volume", "volume_units", "replicate", "condition")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
}

# --- data rows (well, strain, compound, concentration_unit, concentration,
#     compound_volume, media, media_volume, volume_units, replicate, condition) ---
$data = @(
    @("C2", 168, "none", "none", "none", "none", "LB Teknova", 250, "uL", 0, 0),
    @("C3", 168, "none", "none", "none", "none", "LB Teknova", 250, "uL", 1, 0),
    @("C4", 168, "none", "none", "none", "none", "LB Teknova", 250, "uL", 2, 0),
    @("C5", "1x", "none", "none", "none", "none", "LB Teknova", 250, "uL", 0, 1),
    @("C6", "1x", "none", "none", "none", "none", "LB Teknova", 250, "uL", 1, 1),
    @("C7", "1x", "none", "none", "none", "none", "LB Teknova", 250, "uL", 2, 1),
    @("C8", "4x", "none", "none", "none", "none", "LB Teknova", 250, "uL", 0, 2),
    @("C9", "4x", "none", "none", "none", "none", "LB Teknova", 250, "uL", 1, 2),
    @("C10", "4x", "none", "none", "none", "none", "LB Teknova", 250, "uL", 2, 2),
    @("D2", 168, "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 0, 3),
    @("D3", 168, "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 1, 3),
    @("D4", 168, "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 2, 3),
    @("D5", "1x", "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 0, 4),
    @("D6", "1x", "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 1, 4),
    @("D7", "1x", "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 2, 4),
    @("D8", "4x", "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 0, 5),
    @("D9", "4x", "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 1, 5),
    @("D10", "4x", "ampicillin", "ug/mL", 2.5, 1, "LB Teknova", 250, "uL", 2, 5),
    @("E2", 168, "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 0, 6),
    @("E3", 168, "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 1, 6),
    @("E4", 168, "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 2, 6),
    @("E5", "1x", "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 0, 7),
    @("E6", "1x", "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 1, 7),
    @("E7", "1x", "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 2, 7),
    @("E8", "4x", "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 0, 8),
    @("E9", "4x", "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 1, 8),
    @("E10", "4x", "ampicillin", "ug/mL", 5, 2, "LB Teknova", 250, "uL", 2, 8)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# --- cosmetics: column C width, zoom, active/selected sheet -------------
$ws.Columns.Item(3).ColumnWidth = 11.1666667

$ws.Activate()
$excel.ActiveWindow.Zoom = 140
